$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the team header in D11 (this team's roster header) from "DAMU DREAM TEAM" to "JAYAGAN ARMY"
$ws.Range("D11").Value = "JAYAGAN ARMY"

# Fill in Contest 31 (row 43): KKR vs RCB results
$ws.Range("E43").Value = 100
$ws.Range("H43").Value = 40
$ws.Range("K43").Value = 60
$ws.Range("N43").Value = 80
$ws.Range("Q43").Value = 70
$ws.Range("T43").Value = 30
$ws.Range("W43").Value = 50
$ws.Range("Z43").Value = 20
$ws.Range("AC43").Value = 0
